# Create group var for ex3400 switches
#
# The report used to carry duplicate readings for the same three devices
# (rows 2-4 repeated as rows 5-10 with slightly different utilization
# numbers). Trim it down to one reading per device and refresh a few of
# the utilization figures on the rows that remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "master RE cpu/memory utilization" readings for the three
# remaining devices. Assigning with a leading apostrophe keeps these
# text values (matching how the sheet already stores every other number
# in this table as text) instead of letting them get reinterpreted as
# numeric cells; resetting the style back to Normal afterwards avoids
# leaving a stray "quote prefix" cell format behind.
$ws.Range("F2").Value = "'40"
$ws.Range("F2").Style = "Normal"

$ws.Range("G2").Value = "'34"
$ws.Range("G2").Style = "Normal"

$ws.Range("F3").Value = "'69"
$ws.Range("F3").Style = "Normal"

$ws.Range("F4").Value = "'31"
$ws.Range("F4").Style = "Normal"

# The duplicate device rows (5-10) are no longer needed now that each
# device only needs a single row.
$ws.Range("A5:G10").Delete()

# Refresh the selection to cover the surviving data rows now that the
# sheet only runs through row 4.
$ws.Range("A2:G10").Select()
